$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values in the precise order that reproduces the shared-string table
# ordering of the target workbook (first-seen order determines index).
$ws.Range("A12").Value = "CS_CREATE_A_SECURE"
$ws.Range("B12").Value = "Charles Stanley. Create a secure financial future"

$ws.Range("B13").Value = "Don't have an account?,Create one on our website"
$ws.Range("A13").Value = "DON" + [char]0x2019 + "T_HAVE_AN_ACCOUNT_LNK"

$ws.Range("A14").Value = "CS_SECURE_URL"

$ws.Range("A15").Value = "TESAT_DATA"

$ws.Range("B15").Value = "SDFSDFDSFSDFDS"

$ws.Range("B14").Value = [char]0x200E + "charles-stanley.co.uk, secure"

$fc = $ws.Range("A11").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A11:A15"))

$ws.Range("F18").Select()
